$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# Sheet "展览" updates (F column)
$ws1.Range("F4").Value = 3471
$ws1.Range("F5").Value = 3471
$ws1.Range("F6").Value = 243
$ws1.Range("F7").Value = 4998
$ws1.Range("F8").Value = 502
$ws1.Range("F9").Value = 331
$ws1.Range("F10").Value = 190
$ws1.Range("F11").Value = 666
$ws1.Range("F13").Value = 69
$ws1.Range("F15").Value = 684
$ws1.Range("F17").Value = 31
$ws1.Range("F21").Value = 357
$ws1.Range("F22").Value = 4846
$ws1.Range("F23").Value = 41
$ws1.Range("F26").Value = 5968
$ws1.Range("F30").Value = 306
$ws1.Range("F33").Value = 314
$ws1.Range("F34").Value = 110
$ws1.Range("F36").Value = 953
$ws1.Range("F38").Value = 21
$ws1.Range("F40").Value = 841
$ws1.Range("F41").Value = 935
$ws1.Range("F42").Value = 14

# Sheet "全部类型" updates (F column)
$ws4.Range("F8").Value = 3471
$ws4.Range("F9").Value = 3471
$ws4.Range("F10").Value = 243
$ws4.Range("F11").Value = 4998
$ws4.Range("F12").Value = 502
$ws4.Range("F13").Value = 331
$ws4.Range("F14").Value = 190
$ws4.Range("F15").Value = 666
$ws4.Range("F16").Value = 69
$ws4.Range("F18").Value = 684
$ws4.Range("F20").Value = 31
$ws4.Range("F25").Value = 357
$ws4.Range("F26").Value = 4846
$ws4.Range("F27").Value = 41
$ws4.Range("F30").Value = 5968
$ws4.Range("F34").Value = 306
$ws4.Range("F37").Value = 314
$ws4.Range("F39").Value = 110
$ws4.Range("F41").Value = 953
$ws4.Range("F43").Value = 21
$ws4.Range("F45").Value = 841
$ws4.Range("F46").Value = 935
$ws4.Range("F48").Value = 14
